$d = $word.ActiveDocument

# Anchor all subsequent searches to the "Thief" occurrence that is immediately
# followed by ", Bug" and the (soon to move) "_GoBack" bookmark -- this is the
# paragraph the commit touches (the same phrase / words recur later in the
# document for other levels, so we must not just grab the first global hit).
$anchor = $d.Content
$foundAnchor = $anchor.Find.Execute("Thief", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, "", 0)
if (-not $foundAnchor) {
    throw "Could not find the 'Thief' anchor text"
}
$searchStart = $anchor.End
$searchEnd = $d.Content.End

# 1) Remove the word "pop up blocker" from the end of the nearby "Defense
#    types" list, leaving the trailing ", " right after "Nort".
$rng = $d.Range($searchStart, $searchEnd)
$foundPhrase = $rng.Find.Execute("pop up blocker", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if (-not $foundPhrase) {
    throw "Could not find 'pop up blocker' text to remove"
}
$rng.Delete()

# 2) Re-locate the (now trailing) "Nort, " so we can drop the "_GoBack"
#    bookmark right after it (i.e. at the very end of the "Defense types"
#    paragraph).
$rng = $d.Range($searchStart, $d.Content.End)
$foundNort = $rng.Find.Execute("Nort, ", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
if (-not $foundNort) {
    throw "Could not find 'Nort, ' text to anchor the bookmark"
}
$insertPos = $rng.End

# Placing a collapsed bookmark exactly one character before a paragraph mark
# is unreliable in this host, so temporarily insert a sentinel character right
# after the target position, anchor the bookmark next to it (Bookmarks.Add
# re-using the "_GoBack" name also removes the old bookmark automatically),
# then delete the sentinel again.
$sentinelRange = $d.Range($insertPos, $insertPos)
$sentinelRange.InsertAfter("#")

$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$bm = $d.Bookmarks("_GoBack")
$sentinel = $d.Range($bm.End, $bm.End + 1)
$sentinel.Delete()

# 3) Sanity check: "Thief" and ", Bug" must now be contiguous (nothing, not
#    even the bookmark, sits between them any more), and the bookmark must sit
#    at the end of the "Defense types" paragraph.
$check = $d.Range($searchStart - 10, $d.Content.End)
$foundCheck = $check.Find.Execute("Thief, Bug", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "", 0)
if (-not $foundCheck) {
    throw "Expected 'Thief, Bug' to be contiguous after the edit"
}
